$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 94154.75999999999
$ws.Range("I15").Value = 94154.75999999999
$ws.Range("K15").Value = 282464.28
$ws.Range("M15").Value = -282295.28
$ws.Range("H39").Value = 1094.1538
$ws.Range("I39").Value = 1656.7142
$ws.Range("K39").Value = 4970.142599999999
$ws.Range("M39").Value = -4674.142599999999
$ws.Range("H62").Value = 4135996.8
$ws.Range("I62").Value = 6188473
$ws.Range("J62").Value = 31044.334
$ws.Range("K62").Value = 6188473
$ws.Range("L62").Value = 31044.334
$ws.Range("M62").Value = -6187849
$ws.Range("N62").Value = -32292.334
$ws.Range("H65").Value = 4135996.8
$ws.Range("I65").Value = 6188473
$ws.Range("J65").Value = 31044.334
$ws.Range("K65").Value = 30942365
$ws.Range("L65").Value = 155221.67
$ws.Range("M65").Value = -30939245
$ws.Range("N65").Value = -161461.67
$ws.Range("H112").Value = 5246006.5
$ws.Range("J112").Value = 5455823
$ws.Range("L112").Value = 16367469
$ws.Range("N112").Value = -16369685
$ws.Range("H113").Value = 4719.3335
$ws.Range("I113").Value = 2937.7144
$ws.Range("J113").Value = 6278.25
$ws.Range("K113").Value = 2937.7144
$ws.Range("L113").Value = 6278.25
$ws.Range("M113").Value = 316.2856000000002
$ws.Range("N113").Value = -12786.25
$ws.Range("H132").Value = 40927.69
$ws.Range("I132").Value = 55005.684
$ws.Range("K132").Value = 165017.052
$ws.Range("M132").Value = -162487.052
$ws.Range("H133").Value = 44951.5
$ws.Range("J133").Value = 44951.5
$ws.Range("L133").Value = 44951.5
$ws.Range("N133").Value = -55071.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 50000
$ws.Range("J18").Value = 50000
$ws.Range("L18").Value = 50000
$ws.Range("N18").Value = -50644
$ws.Range("H32").Value = 33040.91
$ws.Range("I32").Value = 7171.387
$ws.Range("J32").Value = 90323.42999999999
$ws.Range("K32").Value = 7171.387
$ws.Range("L32").Value = 90323.42999999999
$ws.Range("M32").Value = -6884.387
$ws.Range("N32").Value = -90897.42999999999
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 2966.261
$ws.Range("I61").Value = 1931.2667
$ws.Range("J61").Value = 4906.875
$ws.Range("K61").Value = 1931.2667
$ws.Range("L61").Value = 4906.875
$ws.Range("M61").Value = -1719.2667
$ws.Range("N61").Value = -5330.875
$ws.Range("H81").Value = 44590.5
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 44590.5
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 44590.5
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -46586.5
$ws.Range("H84").Value = 44590.5
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 44590.5
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 133771.5
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -143755.5
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H132").Value = 2400.8223
$ws.Range("I132").Value = 1842.7567
$ws.Range("J132").Value = 4981.875
$ws.Range("K132").Value = 5528.2701
$ws.Range("L132").Value = 14945.625
$ws.Range("M132").Value = -2998.2701
$ws.Range("N132").Value = -20005.625
$ws.Range("H136").Value = 2966.261
$ws.Range("I136").Value = 1931.2667
$ws.Range("J136").Value = 4906.875
$ws.Range("K136").Value = 5793.800099999999
$ws.Range("L136").Value = 14720.625
$ws.Range("M136").Value = -3243.800099999999
$ws.Range("N136").Value = -19820.625

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H26").Value = 26992.334
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H59").Value = 57000
$ws.Range("J59").Value = 57000
$ws.Range("L59").Value = 57000
$ws.Range("N59").Value = -58694
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H134").Value = 3535.7222
$ws.Range("I134").Value = 2711.1304
$ws.Range("J134").Value = 4994.615
$ws.Range("K134").Value = 8133.3912
$ws.Range("L134").Value = 14983.845
$ws.Range("M134").Value = -5598.3912
$ws.Range("N134").Value = -20053.845

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4689.1895
$ws.Range("I31").Value = 1473.421
$ws.Range("K31").Value = 1473.421
$ws.Range("M31").Value = -1178.421
$ws.Range("H34").Value = 4689.1895
$ws.Range("I34").Value = 1473.421
$ws.Range("K34").Value = 1473.421
$ws.Range("M34").Value = -1271.421
$ws.Range("H58").Value = 3024.1
$ws.Range("I58").Value = 1766.5
$ws.Range("K58").Value = 1766.5
$ws.Range("M58").Value = -1563.5
$ws.Range("H68").Value = 22711.428
$ws.Range("J68").Value = 23996.666
$ws.Range("L68").Value = 23996.666
$ws.Range("N68").Value = -25494.666
$ws.Range("H71").Value = 22711.428
$ws.Range("J71").Value = 23996.666
$ws.Range("L71").Value = 71989.99800000001
$ws.Range("N71").Value = -79477.99800000001
$ws.Range("H99").Value = 2204.0322
$ws.Range("I99").Value = 1695.9546
$ws.Range("K99").Value = 1695.9546
$ws.Range("M99").Value = -197.9546
$ws.Range("H122").Value = 2850
$ws.Range("J122").Value = 2850
$ws.Range("L122").Value = 8550
$ws.Range("N122").Value = -13450
$ws.Range("H126").Value = 2204.0322
$ws.Range("I126").Value = 1695.9546
$ws.Range("K126").Value = 5087.8638
$ws.Range("M126").Value = -2617.8638
$ws.Range("H132").Value = 2123.8635
$ws.Range("I132").Value = 1392.4706
$ws.Range("K132").Value = 4177.4118
$ws.Range("M132").Value = -1647.4118
$ws.Range("H134").Value = 2952.4666
$ws.Range("I134").Value = 1735.2273
$ws.Range("J134").Value = 6299.875
$ws.Range("K134").Value = 5205.6819
$ws.Range("L134").Value = 18899.625
$ws.Range("M134").Value = -2670.6819
$ws.Range("N134").Value = -23969.625
$ws.Range("H136").Value = 3024.1
$ws.Range("I136").Value = 1766.5
$ws.Range("K136").Value = 5299.5
$ws.Range("M136").Value = -2749.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5292364.5
$ws.Range("J131").Value = 6290854.5
$ws.Range("L131").Value = 18872563.5
$ws.Range("N131").Value = -18882643.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 33500.332
$ws.Range("I41").Value = 33500.332
$ws.Range("K41").Value = 33500.332
$ws.Range("M41").Value = -33145.332
$ws.Range("H102").Value = 2772.0967
$ws.Range("I102").Value = 2387.4285
$ws.Range("J102").Value = 3579.9
$ws.Range("K102").Value = 2387.4285
$ws.Range("L102").Value = 3579.9
$ws.Range("M102").Value = -765.4285
$ws.Range("N102").Value = -6823.9
$ws.Range("H126").Value = 2815.2173
$ws.Range("I126").Value = 2658.4
$ws.Range("J126").Value = 2858.7778
$ws.Range("K126").Value = 7975.200000000001
$ws.Range("L126").Value = 8576.3334
$ws.Range("M126").Value = -5505.200000000001
$ws.Range("N126").Value = -13516.3334
$ws.Range("H132").Value = 3074.5715
$ws.Range("I132").Value = 2218.8438
$ws.Range("J132").Value = 5812.9
$ws.Range("K132").Value = 6656.5314
$ws.Range("L132").Value = 17438.7
$ws.Range("M132").Value = -4126.5314
$ws.Range("N132").Value = -22498.7
$ws.Range("H139").Value = 44901
$ws.Range("J139").Value = 44901
$ws.Range("L139").Value = 44901
$ws.Range("N139").Value = -55181

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 850
$ws.Range("I22").Value = 775
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 775
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -480
$ws.Range("N22").Value = -1590
$ws.Range("H23").Value = 50000
$ws.Range("I23").Value = 50000
$ws.Range("K23").Value = 50000
$ws.Range("M23").Value = -49770
$ws.Range("H27").Value = 850
$ws.Range("I27").Value = 775
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 775
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -668
$ws.Range("N27").Value = -1214
$ws.Range("H35").Value = 24375
$ws.Range("I35").Value = 31000
$ws.Range("J35").Value = 4500
$ws.Range("K35").Value = 31000
$ws.Range("L35").Value = 4500
$ws.Range("M35").Value = -30664
$ws.Range("N35").Value = -5172
$ws.Range("H122").Value = 3881.1428
$ws.Range("I122").Value = 2952
$ws.Range("J122").Value = 3978.9473
$ws.Range("K122").Value = 8856
$ws.Range("L122").Value = 11936.8419
$ws.Range("M122").Value = -6406
$ws.Range("N122").Value = -16836.8419

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 50000
$ws.Range("I33").Value = 50000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 50000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -49750
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 50000
$ws.Range("I36").Value = 50000
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 50000
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -49750
$ws.Range("N36").ClearContents()
$ws.Range("H122").Value = 1824.8334
$ws.Range("I122").Value = 1527.2222
$ws.Range("J122").Value = 2122.4443
$ws.Range("K122").Value = 4581.6666
$ws.Range("L122").Value = 6367.3329
$ws.Range("M122").Value = -2131.6666
$ws.Range("N122").Value = -11267.3329
$ws.Range("H126").Value = 73978.86
$ws.Range("I126").Value = 85908.664
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 257725.992
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -255255.992
$ws.Range("N126").Value = -12140
